$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Last row of the status table ("Cleantech_10_100kWh" / "SG - CT3"):
# flip the VFT and Dataplicity status cells from online/green to offline/red.
$newColor = 13617910  # RGB(246, 202, 207) -> F6CACF

$vftCell = $t.Cell(7, 2)
$vftCell.Range.Text = "offline"
$vftCell.Shading.BackgroundPatternColor = $newColor

$dataplicityCell = $t.Cell(7, 3)
$dataplicityCell.Range.Text = "offline"
$dataplicityCell.Shading.BackgroundPatternColor = $newColor
